$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link) -- plain text, safe to assign directly
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

# Numeric-looking columns (Price / Volume%) stored as text in the source data;
# force Text number format first so Excel does not auto-convert these to numbers,
# which preserves exact literal formatting (trailing zeros, "%", etc.)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '301.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-3.06%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '35.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-0.12%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.069'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-0.79%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08013'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-2.46%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.934'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-5.98%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '4.042'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-2.08%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '7.758'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-2.21%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9269'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '0.11%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1496'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '32.86%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1896'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.32%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08959'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-2.62%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03456'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-5.55%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09877'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.44%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001397'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-3.18%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005723'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.77%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.533'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.83%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.955'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.02%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3442'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.41%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.64%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.065'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.48%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2397'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '8.66%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04490'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-1.32%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001213'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.14%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004770'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-0.85%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001229'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-1.67%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003021'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-32.11%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01838'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-7.11%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04762'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-3.24%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01061'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '10.27%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.007324'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-4.06%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1329'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-4.07%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002108'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-0.63%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.009728'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-15.98%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006220'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-5.01%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.06%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002098'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.06%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0001998'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.06%'
